# Fill in previously-missing gene names (column B, "genname") for the rows
# whose entrezgene id (column A) already had a lookup result but whose
# name had not been filled in. This mirrors a rerun of the full analysis
# pipeline where a handful of symbol lookups that were previously blank
# are now resolved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    33  = "ZNF503"
    98  = "HS3ST3B1"
    106 = "UBALD2"
    111 = "HAUS1"
    119 = "PPP1R15A"
    121 = "LRATD1"
    123 = "PKDCC"
    142 = "OLIG1"
    143 = "PCP4"
    147 = "NPTXR"
    184 = "CYSTM1"
    185 = "SLC35A4"
    222 = "COL22A1"
    223 = "NRBP2"
    226 = "GABBR2"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
